$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 104, shifting the existing row 104 (and below) down to row 105.
$ws.Rows.Item(104).Insert()

# Copy the style of the date cell from the row above (row 103) to the new row's date cell (D104),
# since column D uses a custom date-time number format (style index 2).
$ws.Range("D103").Copy()
$ws.Range("D104").PasteSpecial(-4122)  # xlPasteFormats

# Populate the newly inserted row 104 with the new record data.
$ws.Range("A104").Value = 5
$ws.Range("B104").Value = "Macroferia Regional de Talca"
$ws.Range("C104").Value = "Maule"
$ws.Range("D104").Value = 44595
$ws.Range("E104").Value = 7
$ws.Range("F104").Value = "Fruta"
$ws.Range("G104").Value = 100108
$ws.Range("H104").Value = "Tropicales y subtropicales"
$ws.Range("I104").Value = 100108002
$ws.Range("J104").Value = "Mango"
$ws.Range("K104").Value = "Sin especificar"
$ws.Range("L104").Value = "Primera"
$ws.Range("M104").Value = 200
$ws.Range("N104").Value = 7000
$ws.Range("O104").Value = 7000
$ws.Range("P104").Value = 7000
$ws.Range("Q104").Value = "`$/bandeja 4 kilos"
$ws.Range("R104").Value = "Perú"
$ws.Range("S104").Value = 1750
$ws.Range("T104").Value = 4
